# Auto-generated Excel COM-interop script
# Updates cached currentAveragePrice / LevePrice / LeveProfit values
# across the ALC, ARM, BSM, CRP, GSM, LTW, and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3585.5454
$ws.Range("I76").Value = 3585.5454
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3585.5454
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3270.5454
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 3585.5454
$ws.Range("I79").Value = 3585.5454
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3585.5454
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2493.5454
$ws.Range("N79").ClearContents()
$ws.Range("H116").Value = 3219.9092
$ws.Range("I116").Value = 2977.375
$ws.Range("J116").Value = 3866.6667
$ws.Range("K116").Value = 2977.375
$ws.Range("L116").Value = 3866.6667
$ws.Range("M116").Value = 464.625
$ws.Range("N116").Value = -10750.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9992.727000000001
$ws.Range("I32").Value = 6012.521
$ws.Range("J32").Value = 31730.77
$ws.Range("K32").Value = 6012.521
$ws.Range("L32").Value = 31730.77
$ws.Range("M32").Value = -5725.521
$ws.Range("N32").Value = -32304.77
$ws.Range("H45").Value = 1939.3
$ws.Range("I45").Value = 1736.625
$ws.Range("K45").Value = 1736.625
$ws.Range("M45").Value = -1359.625
$ws.Range("H63").Value = 2580.25
$ws.Range("I63").Value = 2518.9333
$ws.Range("K63").Value = 2518.9333
$ws.Range("M63").Value = -1832.9333
$ws.Range("H66").Value = 2580.25
$ws.Range("I66").Value = 2518.9333
$ws.Range("K66").Value = 12594.6665
$ws.Range("M66").Value = -9162.666500000001
$ws.Range("H122").Value = 1620.125
$ws.Range("I122").Value = 942.2
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 2826.6
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -376.6000000000004
$ws.Range("N122").Value = -13150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1724.9166
$ws.Range("I105").Value = 1724.9166
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1724.9166
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 22.08339999999998
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2459.6365
$ws.Range("I107").Value = 2982
$ws.Range("J107").Value = 1066.6666
$ws.Range("K107").Value = 2982
$ws.Range("L107").Value = 1066.6666
$ws.Range("M107").Value = -1062
$ws.Range("N107").Value = -4906.6666
$ws.Range("H122").Value = 3900.611
$ws.Range("I122").Value = 4588.875
$ws.Range("J122").Value = 3350
$ws.Range("K122").Value = 13766.625
$ws.Range("L122").Value = 10050
$ws.Range("M122").Value = -11316.625
$ws.Range("N122").Value = -14950
$ws.Range("H127").Value = 32400
$ws.Range("I127").Value = 25000
$ws.Range("J127").Value = 34250
$ws.Range("K127").Value = 25000
$ws.Range("L127").Value = 34250
$ws.Range("M127").Value = -20040
$ws.Range("N127").Value = -44170

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4606.2173
$ws.Range("I70").Value = 4639.737
$ws.Range("J70").Value = 4447
$ws.Range("K70").Value = 4639.737
$ws.Range("L70").Value = 4447
$ws.Range("M70").Value = -4369.737
$ws.Range("N70").Value = -4987
$ws.Range("H73").Value = 4606.2173
$ws.Range("I73").Value = 4639.737
$ws.Range("J73").Value = 4447
$ws.Range("K73").Value = 4639.737
$ws.Range("L73").Value = 4447
$ws.Range("M73").Value = -3703.737
$ws.Range("N73").Value = -6319
$ws.Range("H97").Value = 1255.3846
$ws.Range("I97").Value = 1101.1875
$ws.Range("J97").Value = 1502.1
$ws.Range("K97").Value = 1101.1875
$ws.Range("L97").Value = 1502.1
$ws.Range("M97").Value = -605.1875
$ws.Range("N97").Value = -2494.1
$ws.Range("H102").Value = 1568.9286
$ws.Range("I102").Value = 1172.3158
$ws.Range("J102").Value = 2406.2222
$ws.Range("K102").Value = 1172.3158
$ws.Range("L102").Value = 2406.2222
$ws.Range("M102").Value = 449.6841999999999
$ws.Range("N102").Value = -5650.2222
$ws.Range("H122").Value = 2755.76
$ws.Range("I122").Value = 1467.8334
$ws.Range("J122").Value = 3944.6155
$ws.Range("K122").Value = 4403.5002
$ws.Range("L122").Value = 11833.8465
$ws.Range("M122").Value = -1953.5002
$ws.Range("N122").Value = -16733.8465
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2325
$ws.Range("I40").Value = 2325
$ws.Range("K40").Value = 2325
$ws.Range("M40").Value = -2189
$ws.Range("H61").Value = 2004.1714
$ws.Range("I61").Value = 867.06665
$ws.Range("J61").Value = 2857
$ws.Range("K61").Value = 867.06665
$ws.Range("L61").Value = 2857
$ws.Range("M61").Value = -665.06665
$ws.Range("N61").Value = -3261
$ws.Range("H82").Value = 2916.3333
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 2916.3333
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 2916.3333
$ws.Range("N82").Value = -3638.3333
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 2916.3333
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 2916.3333
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 2916.3333
$ws.Range("N85").Value = -5412.3333
$ws.Range("M85").ClearContents()
$ws.Range("H113").Value = 2004.1714
$ws.Range("I113").Value = 867.06665
$ws.Range("J113").Value = 2857
$ws.Range("K113").Value = 867.06665
$ws.Range("L113").Value = 2857
$ws.Range("M113").Value = 1302.93335
$ws.Range("N113").Value = -7197
$ws.Range("H122").Value = 5278.857
$ws.Range("I122").Value = 5755.4443
$ws.Range("J122").Value = 4421
$ws.Range("K122").Value = 17266.3329
$ws.Range("L122").Value = 13263
$ws.Range("M122").Value = -14816.3329
$ws.Range("N122").Value = -18163

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 23950
$ws.Range("J123").Value = 23950
$ws.Range("L123").Value = 23950
$ws.Range("N123").Value = -33750
$ws.Range("H125").Value = 36510.715
$ws.Range("J125").Value = 36510.715
$ws.Range("L125").Value = 36510.715
$ws.Range("N125").Value = -46350.715

